$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New s_vals data (regenerated to filter save games)
$data = @{
    2 = @(3.286832544864788, 1.655778082260271, 22.3905356188092, 10.19245300693656, 0, 37.52559925287081)
    3 = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 0, 4.358119930609447)
    4 = @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 0, 4.358119930609447)
    5 = @(1.455362044514542, 1.655778082260271, 3.537761648806719, 10.19245300693656, 1, 16.84135478251809)
    6 = @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 1, 5.586269137925634)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}
